$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The work-diary table (Tableau1) gains one more logged entry (row 71).
# Grow the table by one row first so the table ref / autofilter / sheet
# dimension all expand together, just like Excel does when you type into
# the row right below an existing table.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# The previous last entry (row 70, "Bugfix - Tests de la story 'ajout
# d'articles'") was missing its Durée (heures) value - fill it in.
$ws.Range("C70").Value = 1

# New entry: 22.05.2022, Réalisation, 1h, "Fonction de like et dislike"
$ws.Range("A71").Value = 44703
$ws.Range("A71").NumberFormat = "dd/mm/yyyy"
$ws.Range("B71").Value = "Réalisation"
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = "Fonction de like et dislike"

# Leave the selection on the new last cell of the table, as in the saved file.
$ws.Range("F71").Select() | Out-Null
